$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.612.25"
$ws.Range("E2").Value = "  +1.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.411.58"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.59%  "

# Row 5 - BNB
$ws.Range("D5").Value = "567.59"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6 - Solana
$ws.Range("D6").Value = "143.33"
$ws.Range("E6").Value = "  +3.11%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.43%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  -0.55%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.423.80"
$ws.Range("E9").Value = "  +1.15%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +2.37%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.51%  "

# Row 12 - Toncoin
$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  +3.39%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  +3.75%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "26.48"
$ws.Range("E14").Value = "  +2.34%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").Value = "  +3.03%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.855.72"
$ws.Range("E16").Value = "  +0.30%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "61.440.20"
$ws.Range("E17").Value = "  +0.72%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.423.18"
$ws.Range("E18").Value = "  +1.09%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "8.03"
$ws.Range("E19").Value = "  -1.59%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "10.65"
$ws.Range("E20").Value = "  +1.55%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "324.43"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  +1.34%  "

# Row 23 - LEO
$ws.Range("D23").Value = "6.10"
$ws.Range("E23").Value = "  -1.03%  "

# Row 24 - SuiNetwork
$ws.Range("E24").Value = "  +8.12%  "

# Row 25 - Dai
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "65.16"
$ws.Range("E26").Value = "  +1.30%  "

# Row 27 - Bittensor
$ws.Range("D27").Value = "614.94"
$ws.Range("E27").Value = "  +7.04%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").Value = "  +1.47%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0958"
$ws.Range("E29").Value = "  +3.97%  "

# Row 30 - WrappedeETH
$ws.Range("D30").Value = "2.521.02"
$ws.Range("E30").Value = "  -0.43%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  +1.28%  "

# Row 32 - Fetch.AI
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  +4.77%  "

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  +1.16%  "

# Row 35 - swapped to ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.47"
$ws.Range("E35").Value = "  +4.16%  "

# Row 36 - swapped to FirstDigitalUSD
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.55%  "

# Row 37 - Monero
$ws.Range("D37").Value = "153.22"
$ws.Range("E37").Value = "  +0.55%  "

# Row 38 - swapped to NEARProtocol
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.64"
$ws.Range("E38").Value = "  +2.07%  "

# Row 39 - swapped to PolygonEcosystemToken
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.373"
$ws.Range("E39").Value = "  +1.33%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "5.35"
$ws.Range("E40").Value = "  +4.87%  "

# Row 41 - EthereumClassic
$ws.Range("D41").Value = "18.40"
$ws.Range("E41").Value = "  +1.23%  "

# Row 42 - dogwifhat
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +8.65%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +3.62%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  -0.04%  "

# Row 45 - OKB
$ws.Range("D45").Value = "41.83"
$ws.Range("E45").Value = "  +1.60%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -2.74%  "

# Row 47 - Aave
$ws.Range("D47").Value = "143.01"
$ws.Range("E47").Value = "  +0.32%  "

# Row 48 - Filecoin
$ws.Range("D48").Value = "3.55"
$ws.Range("E48").Value = "  +1.11%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "19.95"
$ws.Range("E49").Value = "  +4.09%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "0.595"
$ws.Range("E50").Value = "  +1.97%  "

# Row 51 - Hedera
$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").Value = "  +2.24%  "
